# Applies the cryptos list price/volume update described in the commit
# "Updated cryptos list on Sat Aug 26 17:09:59 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.162.23"
$ws.Range("E2").Value = "  +0.53%  "
# Row 3
$ws.Range("D3").Value = "1.654.83"
$ws.Range("E3").Value = "  +0.31%  "
# Row 4
$ws.Range("E4").Value = "  +0.33%  "
# Row 5
$ws.Range("D5").Value = "'217.74"
$ws.Range("E5").Value = "  -0.10%  "
# Row 6
$ws.Range("D6").Value = "'0.5296"
$ws.Range("E6").Value = "  +1.63%  "
# Row 7
$ws.Range("E7").Value = "  +0.24%  "
# Row 8
$ws.Range("D8").Value = "'0.2625"
$ws.Range("E8").Value = "  +0.45%  "
# Row 9
$ws.Range("D9").Value = "'0.06316"
$ws.Range("E9").Value = "  +0.72%  "
# Row 10
$ws.Range("D10").Value = "'20.38"
$ws.Range("E10").Value = "  -0.44%  "
# Row 11
$ws.Range("D11").Value = "'0.07813"
$ws.Range("E11").Value = "  +1.07%  "
# Row 12
$ws.Range("D12").Value = "'4.514"
$ws.Range("E12").Value = "  +0.93%  "
# Row 13
$ws.Range("D13").Value = "1.661.39"
$ws.Range("E13").Value = "  +0.67%  "
# Row 14
$ws.Range("D14").Value = "1.883.16"
$ws.Range("E14").Value = "  +0.45%  "
# Row 15
$ws.Range("D15").Value = "'0.5484"
$ws.Range("E15").Value = "  +0.91%  "
# Row 16
$ws.Range("D16").Value = "0.0₅8146"
$ws.Range("E16").Value = "  +0.87%  "
# Row 17
$ws.Range("D17").Value = "'65.33"
$ws.Range("E17").Value = "  +1.15%  "
# Row 18
$ws.Range("D18").Value = "26.139.04"
$ws.Range("E18").Value = "  +0.40%  "
# Row 19
$ws.Range("E19").Value = "  +0.21%  "
# Row 20
$ws.Range("D20").Value = "'4.594"
$ws.Range("E20").Value = "  +0.86%  "
# Row 21
$ws.Range("D21").Value = "'190.69"
$ws.Range("E21").Value = "  -0.38%  "
# Row 22
$ws.Range("D22").Value = "'10.07"
$ws.Range("E22").Value = "  +0.39%  "
# Row 23
$ws.Range("D23").Value = "'5.998"
$ws.Range("E23").Value = "  +0.28%  "
# Row 24
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.33%  "
# Row 25
$ws.Range("D25").Value = "'145.14"
$ws.Range("E25").Value = "  +4.65%  "
# Row 26
$ws.Range("D26").Value = "'0.1222"
$ws.Range("E26").Value = "  -0.65%  "
# Row 27
$ws.Range("D27").Value = "'7.202"
$ws.Range("E27").Value = "  -0.48%  "
# Row 28
$ws.Range("D28").Value = "'15.97"
$ws.Range("E28").Value = "  -1.13%  "
# Row 29
$ws.Range("D29").Value = "'1.471"
$ws.Range("E29").Value = "  +5.10%  "
# Row 30
$ws.Range("D30").Value = "'0.05709"
$ws.Range("E30").Value = "  -3.48%  "
# Row 31
$ws.Range("E31").Value = "  +0.03%  "
# Row 32
$ws.Range("D32").Value = "'3.546"
$ws.Range("E32").Value = "  +1.58%  "
# Row 33
$ws.Range("D33").Value = "'3.266"
$ws.Range("E33").Value = "  +1.17%  "
# Row 34
$ws.Range("D34").Value = "'1.586"
$ws.Range("E34").Value = "  +4.73%  "
# Row 35
$ws.Range("D35").Value = "'2.805"
$ws.Range("E35").Value = "  +2.01%  "
# Row 36
$ws.Range("D36").Value = "'2.422"
$ws.Range("E36").Value = "  +0.33%  "
# Row 37
$ws.Range("D37").Value = "'0.9475"
$ws.Range("E37").Value = "  +0.65%  "
# Row 38
$ws.Range("D38").Value = "'0.5724"
$ws.Range("E38").Value = "  +0.69%  "
# Row 39
$ws.Range("E39").Value = "  +0.35%  "
# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8497"
$ws.Range("E40").Value = "  +0.58%  "
# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.788"
$ws.Range("E41").Value = "  -0.96%  "
# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  +0.33%  "
# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'104.04"
$ws.Range("E43").Value = "  +3.40%  "
# Row 44
$ws.Range("D44").Value = "1.038.61"
$ws.Range("E44").Value = "  +3.69%  "
# Row 45
$ws.Range("D45").Value = "1.795.82"
# Row 46
$ws.Range("D46").Value = "'56.66"
$ws.Range("E46").Value = "  +0.18%  "
# Row 47
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.16%  "
# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  -2.96%  "
# Row 49
$ws.Range("D49").Value = "'0.4355"
$ws.Range("E49").Value = "  +1.49%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.857"
$ws.Range("E50").Value = "  +0.06%  "
# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05155"
$ws.Range("E51").Value = "  +0.10%  "
